# Update building block types for the ENA - Raw sequencing reads template.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version number ---
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.0.2"

# --- Annotation table sheet: update header / building block names ---
$ws = $wb.Worksheets.Item("New Table")

# Header row (row 1)
$ws.Range("B1").Value = "Component [Instrument Model]"
$ws.Range("C1").Value = "Term Source REF (NCIT:C177610)"
$ws.Range("D1").Value = "Term Accession Number (NCIT:C177610)"
$ws.Range("H1").Value = "Characteristic [library source]"
$ws.Range("I1").Value = "Term Source REF (GENEPIO:0001965)"
$ws.Range("J1").Value = "Term Accession Number (GENEPIO:0001965)"
$ws.Range("W1").Value = "Output [Data]"

# Data row (row 2) - update ontology term values to match new building blocks
$ws.Range("C2").Value = "OBI"
$ws.Range("D2").Value = "https://bioregistry.io/OBI:0003386"
$ws.Range("H2").Value = "Genomic DNA"
$ws.Range("I2").Value = "NCIT"
$ws.Range("J2").Value = "https://bioregistry.io/NCIT:C95940"
$ws.Range("M2").Value = "https://bioregistry.io/NCIT:C17003"
$ws.Range("P2").Value = "https://bioregistry.io/NCIT:C101294"
$ws.Range("S2").Value = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_0000086"
